$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44242
$ws.Cells.Item(2, 10).Value = 160
$ws.Cells.Item(2, 11).Value = 5000
$ws.Cells.Item(2, 12).Value = 5500
$ws.Cells.Item(2, 13).Value = 5250
$ws.Cells.Item(2, 16).Value = 88

$ws.Cells.Item(3, 4).Value = 44589
$ws.Cells.Item(3, 10).Value = 110
$ws.Cells.Item(3, 11).Value = 5000
$ws.Cells.Item(3, 12).Value = 6000
$ws.Cells.Item(3, 13).Value = 5500
$ws.Cells.Item(3, 16).Value = 92

$ws.Cells.Item(4, 4).Value = 44676
$ws.Cells.Item(4, 10).Value = 120
$ws.Cells.Item(4, 11).Value = 4000
$ws.Cells.Item(4, 12).Value = 4500
$ws.Cells.Item(4, 13).Value = 4250
$ws.Cells.Item(4, 16).Value = 71

$ws.Cells.Item(5, 4).Value = 44760
$ws.Cells.Item(5, 10).Value = 130
$ws.Cells.Item(5, 11).Value = 7000
$ws.Cells.Item(5, 12).Value = 7500
$ws.Cells.Item(5, 13).Value = 7250
$ws.Cells.Item(5, 16).Value = 121

$ws.Cells.Item(6, 4).Value = 44362
$ws.Cells.Item(6, 10).Value = 120
$ws.Cells.Item(6, 11).Value = 8000
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = 8500
$ws.Cells.Item(6, 16).Value = 142

$ws.Cells.Item(7, 4).Value = 44657
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 5000
$ws.Cells.Item(7, 12).Value = 5500
$ws.Cells.Item(7, 13).Value = 5250
$ws.Cells.Item(7, 16).Value = 88

$ws.Cells.Item(8, 4).Value = 45044
$ws.Cells.Item(8, 10).Value = 190
$ws.Cells.Item(8, 11).Value = 4000
$ws.Cells.Item(8, 12).Value = 5000
$ws.Cells.Item(8, 13).Value = 4526
$ws.Cells.Item(8, 16).Value = 75

$ws.Cells.Item(9, 4).Value = 44494
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 13).Value = 5500
$ws.Cells.Item(9, 16).Value = 92

$ws.Cells.Item(10, 4).Value = 44382
$ws.Cells.Item(10, 10).Value = 160
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 13).Value = 7438
$ws.Cells.Item(10, 16).Value = 124

$ws.Cells.Item(11, 4).Value = 44935
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 6000
$ws.Cells.Item(11, 12).Value = 7000
$ws.Cells.Item(11, 13).Value = 6500
$ws.Cells.Item(11, 16).Value = 108

$ws.Cells.Item(12, 4).Value = 44740
$ws.Cells.Item(12, 10).Value = 120
$ws.Cells.Item(12, 11).Value = 6000
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 13).Value = 6500
$ws.Cells.Item(12, 16).Value = 108

$ws.Cells.Item(13, 4).Value = 44281
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 5500
$ws.Cells.Item(13, 12).Value = 6000
$ws.Cells.Item(13, 13).Value = 5750
$ws.Cells.Item(13, 16).Value = 96

$ws.Cells.Item(14, 4).Value = 44785
$ws.Cells.Item(14, 10).Value = 130
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = 7500
$ws.Cells.Item(14, 16).Value = 125

$ws.Cells.Item(15, 4).Value = 44764
$ws.Cells.Item(15, 10).Value = 120
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 8000
$ws.Cells.Item(15, 13).Value = 7500
$ws.Cells.Item(15, 16).Value = 125

$ws.Cells.Item(16, 4).Value = 44963
$ws.Cells.Item(16, 10).Value = 130
$ws.Cells.Item(16, 11).Value = 4000
$ws.Cells.Item(16, 12).Value = 4500
$ws.Cells.Item(16, 13).Value = 4250
$ws.Cells.Item(16, 16).Value = 71

$ws.Cells.Item(17, 4).Value = 44627
$ws.Cells.Item(17, 10).Value = 120
$ws.Cells.Item(17, 11).Value = 4000
$ws.Cells.Item(17, 12).Value = 4500
$ws.Cells.Item(17, 13).Value = 4250
$ws.Cells.Item(17, 16).Value = 71

$ws.Cells.Item(18, 4).Value = 44400
$ws.Cells.Item(18, 10).Value = 120
$ws.Cells.Item(18, 11).Value = 9000
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 13).Value = 9500
$ws.Cells.Item(18, 16).Value = 158

$ws.Cells.Item(20, 4).Value = 44421
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 8000
$ws.Cells.Item(20, 12).Value = 9000
$ws.Cells.Item(20, 13).Value = 8500
$ws.Cells.Item(20, 16).Value = 142

$ws.Cells.Item(21, 4).Value = 44827
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 6000
$ws.Cells.Item(21, 12).Value = 7000
$ws.Cells.Item(21, 13).Value = 6500
$ws.Cells.Item(21, 16).Value = 108

$ws.Cells.Item(22, 4).Value = 44669
$ws.Cells.Item(22, 10).Value = 130
$ws.Cells.Item(22, 11).Value = 4500
$ws.Cells.Item(22, 12).Value = 5000
$ws.Cells.Item(22, 13).Value = 4750
$ws.Cells.Item(22, 16).Value = 79

$ws.Cells.Item(23, 4).Value = 44603
$ws.Cells.Item(23, 10).Value = 140
$ws.Cells.Item(23, 11).Value = 5500
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 13).Value = 5750
$ws.Cells.Item(23, 16).Value = 96

